# Scripts/TC_ACC_10/Default.xlsx rework:
#  - the sample "USERNAME"/"PASSWORD"/"Qatar@2021" login-form scaffolding on the
#    "Global" sheet is torn out
#  - row 2 (the input cells with thin borders + the sample password text) is
#    removed entirely, shrinking the used range back down to A1:B1
#  - the two header cells become plain placeholders "A" / "B"
#  - column widths for A:B collapse to a single, smaller custom width
#  - the view's selection moves to A2 (now an empty cell below the headers)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Drop the whole second row (values + per-cell border styles) and shift
# everything below it up, shrinking the sheet's dimension to A1:B1.
$ws.Range("A2:B2").Delete()

# The former "USERNAME"/"PASSWORD" header labels become plain "A"/"B".
$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "B"

# Columns A and B had individual bestFit widths (10.640625 / 11.49609375);
# they now share one smaller, explicit width.
$ws.Range("A:B").ColumnWidth = 8.6

# Move the sheet's own selection to A2 without disturbing which sheet/tab
# is actually active in the workbook (Action1 stays the active tab).
$previousActive = $wb.ActiveSheet
$ws.Activate()
$ws.Range("A2").Select()
$previousActive.Activate()
